$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All the new rows (9-23) use the same plain "Text" cell style (numFmtId 49)
# as the pre-existing rows 1-8; stamp the format on the whole block up
# front so every new cell - even ones whose value never goes through the
# apostrophe/Style-reset dance below - ends up with that style.
$ws.Range("A9:B23").NumberFormat = "@"

function Set-TextValue($range, $text) {
    # Force a cell to hold a shared string (not a boolean/number) while
    # keeping/ending up with the plain "Text" number format (style s=1),
    # by using a leading apostrophe and then normalising the style back.
    $range.Value = "'" + $text
    $range.Style = "Normal"
    $range.NumberFormat = "@"
}

# Row 15's "B" value is made only of spaces; Excel needs a literal
# quote-prefix (kept, this time) to store it as text, producing the new
# quotePrefix cell style. Do this first so it becomes cellXfs index 2.
$r15 = $ws.Range("B15")
$r15.Value = "'  "
$r15.NumberFormat = "@"

# Rows 2 & 3: swap true/false (lowercase, booleans-as-strings) for
# True/False (capitalised, still text) and move the "CONTAINS...4" input
# up into row 3 (index shifts because of the old true/false removal).
Set-TextValue $ws.Range("B2") "True"
$ws.Range("A3").Value = "CONTAINS:((1|2|3)|4)"
Set-TextValue $ws.Range("B3") "False"

# Rows 9-17: new test rows for the COMP:(...) expressions.
$ws.Range("A9").Value = "COMP:(1|2|a|b|c)"
$ws.Range("B9").Value = "c"

$ws.Range("A10").Value = "COMP:(2|1|a|b|c)"
Set-TextValue $ws.Range("B10") "a"

$ws.Range("A11").Value = "COMP:(1|1|a|b|c)"
Set-TextValue $ws.Range("B11") "b"

$ws.Range("A12").Value = "COMP:(1|2|a|b)"
Set-TextValue $ws.Range("B12") "b"

$ws.Range("A13").Value = "COMP:(2|1|a|b)"
Set-TextValue $ws.Range("B13") "b"

$ws.Range("A14").Value = "COMP:(1|1|a|b)"
Set-TextValue $ws.Range("B14") "a"

$ws.Range("A15").Value = "COMP:(2|1|a)"
# B15 already set above.

$ws.Range("A16").Value = "COMP:(1|1|a)"
Set-TextValue $ws.Range("B16") "a"

$ws.Range("A17").Value = "OR:(IF:(True|False|True)|COMP:(0|4|False|False|True))"
Set-TextValue $ws.Range("B17") "True"

# Rows 18-23: new test rows for RANGE:/REVERSE:/set expressions.
$ws.Range("A18").Value = "(1|2|3|4)#(1|3)"
$ws.Range("B18").Value = "(2|4)"

$ws.Range("A19").Value = "RANGE:(1|7|3)"
$ws.Range("B19").Value = "(1|4|7)"

$ws.Range("A20").Value = "RANGE:(2|5)"
$ws.Range("B20").Value = "(2|3|4|5)"

$ws.Range("A21").Value = "REVERSE:(1|3|5)"
$ws.Range("B21").Value = "(5|3|1)"

$ws.Range("A22").Value = "RANGE:5"
$ws.Range("B22").Value = "(1|2|3|4|5)"

$ws.Range("A23").Value = "REVERSE:RANGE:10#RANGE:(0|8|2)"
$ws.Range("B23").Value = "(10|8|6|4|2)"

$ws.Range("A23").Select() | Out-Null
